# Applies the commit "Atualizado por script em 02-12-2023 20:45":
#  - Rows 116/117 (match Wehen-Kaiserslautern / Magdeburg-Hansa Rostock) had
#    their match-detail columns (F:V) swapped, keeping the leading
#    Indice/pais/torneio/temporada/data_partida columns (A:E) fixed per row.
#  - Rows 131/132 (match Braunschweig-Greuther Furth / Nurnberg-Dusseldorf)
#    had the same F:V swap applied.
#  - A new row 133 was appended with the freshly scraped
#    Magdeburg-Kaiserslautern match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match-detail columns (F:V) between rows 116 and 117 ---
$row116 = $ws.Range("F116:V116").Value2
$row117 = $ws.Range("F117:V117").Value2
$ws.Range("F116:V116").Value = $row117
$ws.Range("F117:V117").Value = $row116

# --- Swap match-detail columns (F:V) between rows 131 and 132 ---
$row131 = $ws.Range("F131:V131").Value2
$row132 = $ws.Range("F132:V132").Value2
$ws.Range("F131:V131").Value = $row132
$ws.Range("F132:V132").Value = $row131

# --- Append new row 133 (Magdeburg vs Kaiserslautern, 2023-12-02) ---
# Copy formatting (bold/border style on col A, datetime format on col E)
# from the row above before filling in the new data.
$ws.Range("A132:V132").Copy()
$ws.Range("A133").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A133").Value = 132
$ws.Range("B133").Value = "germany"
$ws.Range("C133").Value = "2-bundesliga"
$ws.Range("D133").Value = "2023-2024"
$ws.Range("E133").Value = 45262.85416666666
$ws.Range("F133").Value = "Magdeburg"
$ws.Range("G133").Value = 4
$ws.Range("H133").Value = "Kaiserslautern"
$ws.Range("I133").Value = 1
$ws.Range("J133").Value = 2.25
$ws.Range("K133").Value = "26/11/2023 13:42"
$ws.Range("L133").Value = 2
$ws.Range("M133").Value = "02/12/2023 20:27"
$ws.Range("N133").Value = 3.82
$ws.Range("O133").Value = "26/11/2023 13:42"
$ws.Range("P133").Value = 3.86
$ws.Range("Q133").Value = "02/12/2023 20:29"
$ws.Range("R133").Value = 3.01
$ws.Range("S133").Value = "26/11/2023 13:42"
$ws.Range("T133").Value = 3.71
$ws.Range("U133").Value = "02/12/2023 20:25"
$ws.Range("V133").Value = "https://www.betexplorer.com/football/germany/2-bundesliga/magdeburg-kaiserslautern/p8Y2KJPi/"

$wb.Save()
